# Adding the problem size for the LULESH example.
#
# For both the "lulesh" and "lulesh_opt" perfexpert command lines, append
# a new run containing " 45" (same Courier/28pt formatting as the rest of
# the command line). Word's automatic "_GoBack" bookmark (which marks the
# last edited location) moves from its old spot (after the backprop
# example) to the end of the second ("lulesh_opt") line, since that is
# the last place text was typed.

$d = $word.ActiveDocument

function Find-ParagraphByExactText($doc, $exactText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $t2 = $t.Substring(0, $t.Length - 1)   # drop trailing paragraph mark
        } else {
            $t2 = $t
        }
        if ($t2 -eq $exactText) {
            return $p
        }
    }
    return $null
}

# Appends " 45" as its own run (matching the surrounding Courier/28pt
# formatting) right before the end of $para, returning the buffer
# position immediately after the inserted text.
function Add-RunWithSize($doc, $para) {
    $insertPos = $para.Range.End - 1
    $doc.Range($insertPos, $insertPos).InsertAfter(" 45")

    # Force the newly typed text to live in its own <w:r> (rather than
    # being merged back into the preceding run) by toggling a character
    # property on just the new text and back off again.
    $newRange = $doc.Range($insertPos, $insertPos + 3)
    $newRange.Bold = 1
    $newRange.Bold = 0

    return $insertPos + 3
}

# --- "$ OMP_NUM_THREADS=16 perfexpert lulesh" ---
$pLulesh = Find-ParagraphByExactText $d "`$ OMP_NUM_THREADS=16 perfexpert lulesh"
Add-RunWithSize $d $pLulesh | Out-Null

# --- "$ OMP_NUM_THREADS=16 perfexpert lulesh_opt" ---
$pLuleshOpt = Find-ParagraphByExactText $d "`$ OMP_NUM_THREADS=16 perfexpert lulesh_opt"
$endPos = Add-RunWithSize $d $pLuleshOpt

# Move the "_GoBack" bookmark here. A collapsed range sitting exactly on a
# paragraph boundary isn't accepted directly, so insert a temporary
# placeholder character after the target spot, anchor the bookmark next
# to it, then remove the placeholder (the bookmark stays put).
$d.Range($endPos, $endPos).InsertAfter("X")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$target = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
$d.Range($endPos, $endPos + 1).Delete()
